# Updates the "Resumo Inscrições Integrado" sheet with refreshed
# enrollment counts (Inscritos / Pagos / Inscrições homologadas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = 431
$ws.Range("F10").Value = 213
$ws.Range("H10").Value = 213

$ws.Range("E11").Value = 292

$ws.Range("E12").Value = 426

$ws.Range("E14").Value = 107
$ws.Range("F14").Value = 54
$ws.Range("H14").Value = 54

$ws.Range("E16").Value = 177

$ws.Range("E21").Value = 127
$ws.Range("F21").Value = 69
$ws.Range("H21").Value = 69

$ws.Range("E23").Value = 176

$ws.Range("E26").Value = 128

$ws.Range("E27").Value = 296

$ws.Range("E28").Value = 172
$ws.Range("F28").Value = 60
$ws.Range("H28").Value = 60

$ws.Range("E31").Value = 68
$ws.Range("F31").Value = 32
$ws.Range("H31").Value = 32

$ws.Range("E33").Value = 255

$ws.Range("E34").Value = 190
$ws.Range("F34").Value = 116
$ws.Range("H34").Value = 116

$ws.Range("F37").Value = 64
$ws.Range("H37").Value = 64

$ws.Range("E39").Value = 167
$ws.Range("F39").Value = 78
$ws.Range("H39").Value = 78

$ws.Range("E41").Value = 348

$ws.Range("E42").Value = 319
$ws.Range("F42").Value = 168
$ws.Range("H42").Value = 168

$ws.Range("E43").Value = 102

$ws.Range("E44").Value = 273

$ws.Range("E45").Value = 126

$ws.Range("E46").Value = 275

$ws.Range("E47").Value = 389
$ws.Range("F47").Value = 191
$ws.Range("H47").Value = 191

$ws.Range("E48").Value = 182

$ws.Range("E50").Value = 225
